$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Fill in the previously-blank columns B-K and M on row 20 with "nan"
# (matches the other historical log rows in this sheet)
foreach ($col in @("B","C","D","E","F","G","H","I","J","K")) {
    $ws.Range($col + "20").Value = "nan"
}
$ws.Range("M20").Value = "nan"

# Add the new service-log entry in row 21
# (use a text formula + paste-values so "9" is stored as text, not a number,
# without leaving a stray NumberFormat style behind)
$a21 = $ws.Range("A21")
$a21.Formula = "=""9"""
$a21.Copy()
$a21.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("L21").Value = "23/11/2025"
$ws.Range("M21").Value = "قطع سير كويلر مسنن دبل 700"
$ws.Range("N21").Value = "تم تغير سير  دوبل700(محمد نعيم)"
$ws.Range("O21").Value = "فني"
